{"js": "// Change the highlight color of the \"All communication between scenes...\"\n// requirement bullet from yellow to green (marking it complete).\nconst body = context.document.body;\n\n// Locate the bullet by a snippet of its (unique) text.\nconst results = body.search(\"All communication between scenes should be done\", { matchCase: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph not found.\");\n}\n\n// Expand to the whole paragraph so every run (and the paragraph mark\n// itself) gets the new highlight, matching the original yellow coverage.\nconst paragraph = results.items[0].paragraphs.getFirst();\nconst range = paragraph.getRange();\n\n// Pure green maps to OOXML's <w:highlight w:val=\"green\"/> (as opposed to\n// named colors like \"Green\" which map to the darker \"darkGreen\" swatch).\nrange.font.highlightColor = \"#00FF00\";\n\nawait context.sync();\n", "ps1": "# Change the highlight color of the \"All communication between scenes...\"\n# requirement bullet from yellow to green (marking it complete).\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"All communication between scenes should be done\")\n\nif (-not $found) {\n    throw \"Target paragraph not found.\"\n}\n\n# Grab the whole paragraph (this also covers the trailing paragraph-mark\n# character) so every run, as well as the paragraph mark itself, gets the\n# new highlight - matching the original yellow coverage.\n$paragraph = $searchRange.Paragraphs(1)\n$paragraphRange = $paragraph.Range\n\n# wdBrightGreen (4) maps to OOXML's <w:highlight w:val=\"green\"/>.\n$paragraphRange.Font.HighlightColorIndex = 4\n"}
